# Update Fgf13-Scn8a LR-pair data with new TPM values, adding the "ECs" sending/target cluster
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf13"
$ws.Range("C2").Value = "Scn8a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05166833333333334
$ws.Range("H2").Value = 0.155005
$ws.Range("I2").Value = 0.01629802083219426
$ws.Range("J2").Value = 0.01629802083219426
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.018839
$ws.Range("N2").Value = 0.056517
$ws.Range("O2").Value = 0.01572277703685034
$ws.Range("P2").Value = 0.01572277703685034
$ws.Range("Q2").Value = 0.0009733797316666666
$ws.Range("R2").Value = 0.008760417585
$ws.Range("S2").Value = 0.0002562501476865324
$ws.Range("T2").Value = 0.0002562501476865324

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf13"
$ws.Range("C3").Value = "Scn8a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05166833333333334
$ws.Range("H3").Value = 0.155005
$ws.Range("I3").Value = 0.01629802083219426
$ws.Range("J3").Value = 0.01629802083219426
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8653686666666666
$ws.Range("N3").Value = 2.596106
$ws.Range("O3").Value = 0.7222250969094145
$ws.Range("P3").Value = 0.7222250969094145
$ws.Range("Q3").Value = 0.04471215672555556
$ws.Range("R3").Value = 0.40240941053
$ws.Range("S3").Value = 0.01177083967496316
$ws.Range("T3").Value = 0.01177083967496316

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf13"
$ws.Range("C4").Value = "Scn8a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05166833333333334
$ws.Range("H4").Value = 0.155005
$ws.Range("I4").Value = 0.01629802083219426
$ws.Range("J4").Value = 0.01629802083219426
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3139903333333333
$ws.Range("N4").Value = 0.941971
$ws.Range("O4").Value = 0.2620521260537351
$ws.Range("P4").Value = 0.2620521260537351
$ws.Range("Q4").Value = 0.01622335720611111
$ws.Range("R4").Value = 0.146010214855
$ws.Range("S4").Value = 0.004270931009544571
$ws.Range("T4").Value = 0.004270931009544571

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf13"
$ws.Range("C5").Value = "Scn8a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03538133333333333
$ws.Range("H5").Value = 0.106144
$ws.Range("I5").Value = 0.01116052464896247
$ws.Range("J5").Value = 0.01116052464896247
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.018839
$ws.Range("N5").Value = 0.056517
$ws.Range("O5").Value = 0.01572277703685034
$ws.Range("P5").Value = 0.01572277703685034
$ws.Range("Q5").Value = 0.0006665489386666666
$ws.Range("R5").Value = 0.005998940448
$ws.Range("S5").Value = 0.0001754744406699093
$ws.Range("T5").Value = 0.0001754744406699093

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf13"
$ws.Range("C6").Value = "Scn8a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03538133333333333
$ws.Range("H6").Value = 0.106144
$ws.Range("I6").Value = 0.01116052464896247
$ws.Range("J6").Value = 0.01116052464896247
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8653686666666666
$ws.Range("N6").Value = 2.596106
$ws.Range("O6").Value = 0.7222250969094145
$ws.Range("P6").Value = 0.7222250969094145
$ws.Range("Q6").Value = 0.03061789725155555
$ws.Range("R6").Value = 0.275561075264
$ws.Range("S6").Value = 0.008060410996156828
$ws.Range("T6").Value = 0.00806041099615683

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf13"
$ws.Range("C7").Value = "Scn8a"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03538133333333333
$ws.Range("H7").Value = 0.106144
$ws.Range("I7").Value = 0.01116052464896247
$ws.Range("J7").Value = 0.01116052464896247
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3139903333333333
$ws.Range("N7").Value = 0.941971
$ws.Range("O7").Value = 0.2620521260537351
$ws.Range("P7").Value = 0.2620521260537351
$ws.Range("Q7").Value = 0.01110939664711111
$ws.Range("R7").Value = 0.099984569824
$ws.Range("S7").Value = 0.002924639212135731
$ws.Range("T7").Value = 0.002924639212135731

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf13"
$ws.Range("C8").Value = "Scn8a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.083171666666667
$ws.Range("H8").Value = 9.249515
$ws.Range("I8").Value = 0.9725414545188432
$ws.Range("J8").Value = 0.9725414545188433
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.018839
$ws.Range("N8").Value = 0.056517
$ws.Range("O8").Value = 0.01572277703685034
$ws.Range("P8").Value = 0.01572277703685034
$ws.Range("Q8").Value = 0.05808387102833333
$ws.Range("R8").Value = 0.522754839255
$ws.Range("S8").Value = 0.01529105244849389
$ws.Range("T8").Value = 0.01529105244849389

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf13"
$ws.Range("C9").Value = "Scn8a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.083171666666667
$ws.Range("H9").Value = 9.249515
$ws.Range("I9").Value = 0.9725414545188432
$ws.Range("J9").Value = 0.9725414545188433
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8653686666666666
$ws.Range("N9").Value = 2.596106
$ws.Range("O9").Value = 0.7222250969094145
$ws.Range("P9").Value = 0.7222250969094145
$ws.Range("Q9").Value = 2.668080154287777
$ws.Range("R9").Value = 24.01272138859
$ws.Range("S9").Value = 0.7023938462382945
$ws.Range("T9").Value = 0.7023938462382946

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf13"
$ws.Range("C10").Value = "Scn8a"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.083171666666667
$ws.Range("H10").Value = 9.249515
$ws.Range("I10").Value = 0.9725414545188432
$ws.Range("J10").Value = 0.9725414545188433
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3139903333333333
$ws.Range("N10").Value = 0.941971
$ws.Range("O10").Value = 0.2620521260537351
$ws.Range("P10").Value = 0.2620521260537351
$ws.Range("Q10").Value = 0.9680860993405556
$ws.Range("R10").Value = 8.712774894065001
$ws.Range("S10").Value = 0.2548565558320548
$ws.Range("T10").Value = 0.2548565558320549

